$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.516.07'
$ws.Range('E2').Value = '  +2.28%  '
$ws.Range('D3').Value = '3.562.39'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''598.09'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = '''172.11'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.32%  '
$ws.Range('D7').Value = '3.556.05'
$ws.Range('E7').Value = '  +1.30%  '
$ws.Range('D8').Value = '''0.613'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('E10').Value = '  +4.63%  '
$ws.Range('D11').Value = '''7.39'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +9.27%  '
$ws.Range('E12').Value = '  +1.75%  '
$ws.Range('D13').Value = '''46.34'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').Value = '4.135.28'
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('D16').Value = '''8.36'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.63%  '
$ws.Range('D17').Value = '''610.69'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '3.565.11'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('D19').Value = '70.534.48'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').Value = '''17.34'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = '''0.881'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').Value = '''9.22'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -16.87%  '
$ws.Range('D24').Value = '''15.75'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('E25').Value = '  +0.33%  '
$ws.Range('D26').Value = '''3.74'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.34%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''2.61'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '''33.83'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.02%  '
$ws.Range('D30').Value = '''9.05'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('D31').Value = '''8.26'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.79%  '
$ws.Range('D32').Value = '''3.04'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.37%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '''661.50'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +7.80%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''7.10'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +3.04%  '
$ws.Range('D35').Value = '''1.30'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.37%  '
$ws.Range('D36').Value = '''3.61'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +4.67%  '
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('D38').Value = '''10.77'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.54%  '
$ws.Range('E39').Value = '  +7.49%  '
$ws.Range('D40').Value = '''57.38'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.50%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('E42').Value = '  +4.14%  '
$ws.Range('D43').Value = '3.380.71'
$ws.Range('E43').Value = '  +0.02%  '
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('D45').Value = '0.0₃0704'
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('D46').Value = '''32.77'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.35%  '
$ws.Range('D47').Value = '''2.94'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +7.43%  '
$ws.Range('D48').Value = '''2.63'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +4.66%  '
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').Value = '''132.28'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('E51').Value = '  -0.10%  '
